$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wednesday column block (M-P): trip header/details for the Metcalfe's Market run
$ws.Range("N34").Value = "4:45 AM MEET AT THE OFFICE"
$ws.Range("N35").Value = "6:00 AM START"
$ws.Range("N36").Value = "DC5-FINANCIAL"
$ws.Range("N37").Value = "METCALFE'S MARKET #2600, WEST MADISON - LIFO"
$ws.Range("N38").Value = "7455 MINERAL POINT RD"
$ws.Range("N39").Value = "https://goo.gl/maps/NW4tsxQQBWF2"

# Wednesday staff list rows 41-55
$ws.Range("M41").Value = "1)"
$ws.Range("N41").Value = "Sarah"
$ws.Range("O41").Value = "@ Store, Equip"

$ws.Range("M42").Value = "2)"
$ws.Range("N42").Value = "Katherine"

$ws.Range("M43").Value = "3)"
$ws.Range("N43").Value = "Lashaun"
$ws.Range("O43").Value = "Driver,`nPrius"

$ws.Range("M44").Value = "4)"
$ws.Range("N44").Value = "Angela"
$ws.Range("O44").Value = "@ Store"

$ws.Range("M45").Value = "5)"
$ws.Range("N45").Value = "Anisha"
$ws.Range("O45").Value = "@ Store"

$ws.Range("M46").Value = "6)"
$ws.Range("N46").Value = "Ashley P"
$ws.Range("O46").Value = "@ Store"

$ws.Range("M47").Value = "7)"
$ws.Range("N47").Value = "Eva"
$ws.Range("O47").Value = "@ Store"

$ws.Range("M48").Value = "8)"
$ws.Range("N48").Value = "Evelin"
$ws.Range("O48").Value = "@ Store,`nAfter Liq Store"

$ws.Range("M49").Value = "9)"
$ws.Range("N49").Value = "Joseph"
$ws.Range("O49").Value = "@ Store"

$ws.Range("M50").Value = "10)"
$ws.Range("N50").Value = "Lori"
$ws.Range("O50").Value = "@ Store"

$ws.Range("M51").Value = "11)"
$ws.Range("N51").Value = "Michael N"
$ws.Range("O51").Value = "@ Store"

$ws.Range("M52").Value = "12)"
$ws.Range("N52").Value = "Nate"
$ws.Range("O52").Value = "@ Store,`nAfter Liq Store"

$ws.Range("M53").Value = "13)"
$ws.Range("N53").Value = "Qiana"
$ws.Range("O53").Value = "@ Store"

$ws.Range("M54").Value = "14)"
$ws.Range("N54").Value = "Savannah"
$ws.Range("O54").Value = "@ Store"

$ws.Range("M55").Value = "15)"
$ws.Range("N55").Value = "Spencer P"
$ws.Range("O55").Value = "@ Store"
